$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "Global Peace Index 2023"
$ws.Range("D12").Value = "Poverty `$6.85 a day (% of population)"
$ws.Range("D16").Value = "doing business score"
$ws.Range("D18").Value = "Armed forces (% labor force)"
$ws.Range("D23").Value = "Corruption index 2023"
$ws.Range("D24").Value = "Poverty `$3.65 a day (% of population)"
$ws.Range("D26").Value = "Poverty `$2.15 a day (% of population)"
$ws.Range("D35").Value = "Unemployment advanced education"
$ws.Range("D37").Value = "Time to start a business (days)"
$ws.Range("D48").Value = "GNI per capita (2015 US$)"
